$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 14: "Πως CertiKOS προστατεύει ;" -> "Πως το CertiKOS προστατεύει ;"
# (Title shape has a single run, safe to set the whole TextRange text.)
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$s14Title = $s14.Shapes.Item(1)
$s14Title.TextFrame.TextRange.Text = "Πως το CertiKOS προστατεύει ;"

# Slide 14: "...απεριόριστους πόους..." -> "...απεριόριστους πόρους..."
# (Second paragraph, first run of the content placeholder.)
$s14Body = $s14.Shapes.Item(2)
$s14Run = $s14Body.TextFrame.TextRange.Paragraphs(2, 1).Runs(1, 1)
$s14Run.Text = "O Hacker δεν μπορεί να δεσμεύσει απεριόριστους πόρους, καθώς η "

# ---------------------------------------------------------------------------
# Slide 5: "...πολυνηματική λειτορυγία..." -> "...πολυνηματική λειτουργία..."
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5Body = $s5.Shapes.Item(2)
$s5Run = $s5Body.TextFrame.TextRange.Paragraphs(3, 1).Runs(1, 1)
$s5Run.Text = "Περιλαμβάει 6500 γραμμές C και x86 assembly, και η απόδειξη λειτουργικής ορθότητας για πολυνηματική λειτουργία ολοκληρώθηκε σε λιγότερο από 2 ανθρωποέτη. Είναι η πρώτη τέτοια απόδειξη για γενικού σκοπού OS πυρήνα."

# ---------------------------------------------------------------------------
# Slide 7: "...ενός εραλείου λογικής..." -> "...ενός εργαλείου λογικής..."
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7Body = $s7.Shapes.Item(2)
$s7Run = $s7Body.TextFrame.TextRange.Paragraphs(10, 1).Runs(5, 1)
$s7Run.Text = ", ενός εργαλείου λογικής και μαθηματικής απόδειξης, για να εξασφαλιστεί η ακεραιότητα του μικροπυρήνα. "

# ---------------------------------------------------------------------------
# Slide 8: "Υπερόπτες" -> "Hypervisors"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8Body = $s8.Shapes.Item(2)
$s8Run = $s8Body.TextFrame.TextRange.Paragraphs(2, 1).Runs(1, 1)
$s8Run.Text = "Hypervisors"

# ---------------------------------------------------------------------------
# Slide 17: center-align the title paragraph ("σύνδεσμος για τον αντίστοιχο
# κώδικα στο github") - text itself is unchanged.
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$s17Title = $s17.Shapes.Item(1)
$s17Title.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# Slide 17: "link " -> "link https://github.com/georgiabasa/project_security"
$s17Body = $s17.Shapes.Item(2)
$s17Run1 = $s17Body.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$s17Run1.Text = "link https://github.com/georgiabasa/project_security"

# Slide 17: "workspace/certikos/kernel" -> "workspace/certikos/kernel."
$s17Run2 = $s17Body.TextFrame.TextRange.Paragraphs(3, 1).Runs(2, 1)
$s17Run2.Text = "workspace/certikos/kernel."

# ---------------------------------------------------------------------------
# Theme: recolor the "Office" colour scheme into the "Calligraphy" palette.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$clrScheme = $master.Theme.ThemeColorScheme
$clrScheme.Colors(3).RGB = 70721        # dk2      411401
$clrScheme.Colors(4).RGB = 15132415     # lt2      FFE6E6
$clrScheme.Colors(5).RGB = 4737698      # accent1  A24A48
$clrScheme.Colors(6).RGB = 6067122      # accent2  B2935C
$clrScheme.Colors(7).RGB = 10132074     # accent3  6A9A9A
$clrScheme.Colors(8).RGB = 8894386      # accent4  B2B787
$clrScheme.Colors(9).RGB = 4940945      # accent5  91644B
$clrScheme.Colors(10).RGB = 7752293     # accent6  654A76
$clrScheme.Colors(11).RGB = 43008       # hlink    00A800
$clrScheme.Colors(12).RGB = 16711935    # folHlink FF00FF
